# Auto update Excel log
# Appends newly-collected sensor readings to the "Humidity" and
# "Temperature" worksheets (the log keeps growing from where it left off).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Humidity sheet: append rows 203-214
# ---------------------------------------------------------------------------
$wsHumidity = $wb.Worksheets.Item("Humidity")

$humidityData = @()
$humidityData += ,@("2026-02-01","14:18:15","14:00","Bathroom","76.6%","Active")
$humidityData += ,@("2026-02-01","14:18:18","14:00","Bathroom","77.6%","Active")
$humidityData += ,@("2026-02-01","14:18:23","14:00","Bathroom","76.6%","Active")
$humidityData += ,@("2026-02-01","14:18:28","14:00","Bathroom","77.5%","Active")
$humidityData += ,@("2026-02-01","14:18:33","14:00","Bathroom","76.7%","Active")
$humidityData += ,@("2026-02-01","14:18:43","14:00","Bathroom","76.7%","Active")
$humidityData += ,@("2026-02-01","14:18:48","14:00","Bathroom","77.6%","Active")
$humidityData += ,@("2026-02-01","14:18:53","14:00","Bathroom","76.7%","Active")
$humidityData += ,@("2026-02-01","14:18:58","14:00","Bathroom","77.6%","Active")
$humidityData += ,@("2026-02-01","14:19:03","14:00","Bathroom","76.7%","Active")
$humidityData += ,@("2026-02-01","14:19:08","14:00","Bathroom","77.6%","Active")
$humidityData += ,@("2026-02-01","14:19:13","14:00","Bathroom","76.6%","Active")

$humidityStartRow = 203
$humidityEndRow = $humidityStartRow + $humidityData.Count - 1
$humidityRange = $wsHumidity.Range("A" + $humidityStartRow + ":F" + $humidityEndRow)

# Force text storage so date-like / percentage-like strings ("2026-02-01",
# "76.6%") are kept as literal text instead of being reinterpreted as a date
# serial number or a numeric percentage.
$humidityRange.NumberFormat = "@"

for ($i = 0; $i -lt $humidityData.Count; $i++) {
    $row = $humidityStartRow + $i
    $values = $humidityData[$i]
    $wsHumidity.Cells.Item($row, 1).Value = $values[0]
    $wsHumidity.Cells.Item($row, 2).Value = $values[1]
    $wsHumidity.Cells.Item($row, 3).Value = $values[2]
    $wsHumidity.Cells.Item($row, 4).Value = $values[3]
    $wsHumidity.Cells.Item($row, 5).Value = $values[4]
    $wsHumidity.Cells.Item($row, 6).Value = $values[5]
}

# Drop the explicit text format again so the new cells fall back to the
# workbook's default (unstyled) cell format, matching the rest of the log.
$humidityRange.ClearFormats()

# ---------------------------------------------------------------------------
# Temperature sheet: append rows 123-135
# ---------------------------------------------------------------------------
$wsTemperature = $wb.Worksheets.Item("Temperature")

$temperatureData = @()
$temperatureData += ,@("2026-02-01","14:18:15","14:00","Bathroom","29.5C","Active")
$temperatureData += ,@("2026-02-01","14:18:16","14:00","Bathroom","29.5C","Active")
$temperatureData += ,@("2026-02-01","14:18:18","14:00","Bathroom","29.5C","Active")
$temperatureData += ,@("2026-02-01","14:18:23","14:00","Bathroom","29.5C","Active")
$temperatureData += ,@("2026-02-01","14:18:29","14:00","Bathroom","29.4C","Active")
$temperatureData += ,@("2026-02-01","14:18:34","14:00","Bathroom","29.5C","Active")
$temperatureData += ,@("2026-02-01","14:18:44","14:00","Bathroom","29.4C","Active")
$temperatureData += ,@("2026-02-01","14:18:49","14:00","Bathroom","29.4C","Active")
$temperatureData += ,@("2026-02-01","14:18:54","14:00","Bathroom","29.5C","Active")
$temperatureData += ,@("2026-02-01","14:18:59","14:00","Bathroom","29.5C","Active")
$temperatureData += ,@("2026-02-01","14:19:04","14:00","Bathroom","29.5C","Active")
$temperatureData += ,@("2026-02-01","14:19:09","14:00","Bathroom","29.4C","Active")
$temperatureData += ,@("2026-02-01","14:19:14","14:00","Bathroom","29.4C","Active")

$temperatureStartRow = 123
$temperatureEndRow = $temperatureStartRow + $temperatureData.Count - 1
$temperatureRange = $wsTemperature.Range("A" + $temperatureStartRow + ":F" + $temperatureEndRow)

$temperatureRange.NumberFormat = "@"

for ($i = 0; $i -lt $temperatureData.Count; $i++) {
    $row = $temperatureStartRow + $i
    $values = $temperatureData[$i]
    $wsTemperature.Cells.Item($row, 1).Value = $values[0]
    $wsTemperature.Cells.Item($row, 2).Value = $values[1]
    $wsTemperature.Cells.Item($row, 3).Value = $values[2]
    $wsTemperature.Cells.Item($row, 4).Value = $values[3]
    $wsTemperature.Cells.Item($row, 5).Value = $values[4]
    $wsTemperature.Cells.Item($row, 6).Value = $values[5]
}

$temperatureRange.ClearFormats()
